$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that hold numeric-looking text (prices) must be forced to Text
# format first, otherwise Excel auto-converts them to numbers and mangles
# formatting (e.g. drops trailing zeros, introduces float artifacts).

$ws.Range("D2").Value = '57.595.19'
$ws.Range("E2").Value = '  +1.89%  '
$ws.Range("D3").Value = '3.011.34'
$ws.Range("E3").Value = '  +0.38%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '510.01'
$ws.Range("E5").Value = '  +0.56%  '
$ws.Range("E6").Value = '  +1.76%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.438'
$ws.Range("E8").Value = '  +1.41%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.53'
$ws.Range("E9").Value = '  -0.57%  '
$ws.Range("E10").Value = '  +1.71%  '
$ws.Range("E11").Value = '  +4.16%  '
$ws.Range("D12").Value = '3.522.11'
$ws.Range("E12").Value = '  +0.19%  '
$ws.Range("E13").Value = '  +0.85%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.44'
$ws.Range("E14").Value = '  +3.94%  '
$ws.Range("E15").Value = '  +6.79%  '
$ws.Range("D16").Value = '57.568.28'
$ws.Range("E16").Value = '  +1.78%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.22'
$ws.Range("E17").Value = '  +7.69%  '
$ws.Range("D18").Value = '3.006.01'
$ws.Range("E18").Value = '  +0.07%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.83'
$ws.Range("E19").Value = '  +3.37%  '
$ws.Range("E20").Value = '  +2.08%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '331.36'
$ws.Range("E21").Value = '  +1.54%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.994'
$ws.Range("E22").Value = '  -0.30%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.499'
$ws.Range("E23").Value = '  +4.44%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.50'
$ws.Range("E24").Value = '  +3.42%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.169'
$ws.Range("E25").Value = '  +0.08%  '
$ws.Range("E26").Value = '  -0.28%  '
$ws.Range("D27").Value = '0.0₃0923'
$ws.Range("E27").Value = '  +1.24%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.79'
$ws.Range("E28").Value = '  +4.15%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.33'
$ws.Range("E29").Value = '  +4.21%  '
$ws.Range("E30").Value = '  +2.44%  '
$ws.Range("E31").Value = '  -5.15%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.59'
$ws.Range("E32").Value = '  -0.09%  '
$ws.Range("E33").Value = '  +5.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '153.54'
$ws.Range("E34").Value = '  -1.29%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.86'
$ws.Range("E35").Value = '  +4.72%  '
$ws.Range("E36").Value = '  +1.41%  '
$ws.Range("B37").Value = 'EnergySwap'
$ws.Range("C37").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '24.46'
$ws.Range("E37").Value = '  +2.23%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0683'
$ws.Range("E38").Value = '  +1.26%  '
$ws.Range("D39").Value = '3.042.57'
$ws.Range("E39").Value = '  +0.28%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '37.32'
$ws.Range("E40").Value = '  +1.81%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.84'
$ws.Range("E41").Value = '  +6.76%  '
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.650'
$ws.Range("E43").Value = '  +0.06%  '
$ws.Range("D44").Value = '2.235.47'
$ws.Range("E44").Value = '  -1.20%  '
$ws.Range("E45").Value = '  +0.51%  '
$ws.Range("E46").Value = '  -1.00%  '
$ws.Range("E47").Value = '  +4.90%  '
$ws.Range("E48").Value = '  +2.03%  '
$ws.Range("E49").Value = '  +2.18%  '
$ws.Range("E50").Value = '  -6.70%  '
$ws.Range("E51").Value = '  +2.67%  '
